$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.772.10'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.918.93'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.01'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4927'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3009'
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06792'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.921.78'
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.34'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07349'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.228'
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.18'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6807'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.751.80'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008018'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.69'
$ws.Range('E18').Value = '  +4.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.182.66'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.428'
$ws.Range('E21').Value = '  +12.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '201.41'
$ws.Range('E23').Value = '  +8.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.363'
$ws.Range('E24').Value = '  +4.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.712'
$ws.Range('E25').Value = '  +3.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.36'
$ws.Range('E26').Value = '  +3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.93'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.981'
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.478'
$ws.Range('E29').Value = '  +5.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.372'
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09197'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.099'
$ws.Range('E32').Value = '  +2.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05345'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7512'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.132'
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.703'
$ws.Range('E36').Value = '  -1.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01870'
$ws.Range('E37').Value = '  +1.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.731'
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9339'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.102'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4534'
$ws.Range('E41').Value = '  +2.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.03'
$ws.Range('E42').Value = '  +26.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.999'
$ws.Range('E43').Value = '  +4.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '107.89'
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1407'
$ws.Range('E45').Value = '  +5.30%  '
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.772'
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '36.12'
$ws.Range('E48').Value = '  +7.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.165'
$ws.Range('E49').Value = '  +5.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05929'
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4075'
$ws.Range('E51').Value = '  +3.84%  '
